# Duplicate the "52" sheet into a new "53" sheet (LQR sim file change),
# placing it at the end of the workbook, updating its filenumber value,
# and making it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("52")
$sourceSheet.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "53"
$newSheet.Range("B1").Value = 53
$newSheet.Activate()
